# Creacion de script, para agregar eficiencia
# Adds a new column "Eficiencia_Terminal" (I) to the estudiantes_graduados
# sheet. For each student row, flags (1/0) whether the student graduated
# within the "normal" 5-year duration of the program, using the
# año_ingreso (E), fecha_egreso (F) and the entry/exit semester
# (termino_ingreso G / termino_egreso H) to resolve the edge case where
# exactly 6 calendar years elapsed but the student still finished within
# 10 semesters.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 124

# Header cell: reuse the same look (bold, centered, thin border) as the
# other header cells by copying the format from H1, then set the text.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Value = "Eficiencia_Terminal"

for ($r = 2; $r -le $lastRow; $r++) {
    $anioIngreso = [int]$ws.Cells.Item($r, 5).Value()
    $anioEgreso = [int]$ws.Cells.Item($r, 6).Value()
    $terminoIngreso = $ws.Cells.Item($r, 7).Value()
    $terminoEgreso = $ws.Cells.Item($r, 8).Value()

    $duracionAnios = $anioEgreso - $anioIngreso

    $eficiencia = 0
    if ($duracionAnios -eq 5) {
        # Exactly 5 years -> always within the normal duration.
        $eficiencia = 1
    } elseif ($duracionAnios -eq 6) {
        # 6 years: only "on time" (<=10 semesters) when entry/exit terms
        # differ (i.e. the student crossed only one extra half-year).
        if ($terminoIngreso -ne $terminoEgreso) {
            $eficiencia = 1
        } else {
            $eficiencia = 0
        }
    } else {
        $eficiencia = 0
    }

    $ws.Cells.Item($r, 9).Value = $eficiencia
}
